$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update description text for 2007.113.1 (row 122, column B)
$ws.Range("B122").Value = "Painting, <em>Once Upon a Time in the West</em>, circa 1968.`nThis painting is the original artwork created for Paramount Picture's Swedish movie poster ""Harmonica En Hamnare."""
$ws.Rows.Item(122).RowHeight = 12.75

# Update description text for 2008.48.1 (row 128, column B)
$ws.Range("B128").Value = "Painting by Franz Arthur Bischoff, <em>The Docks at San Pedro</em>, probably circa 1900."
$ws.Rows.Item(128).RowHeight = 12.75

# Update description text for 2012.37.2 (row 142, column B) - fix "Course" to "Couse"
$ws.Range("B142").Value = "Painting by Eanger Irving Couse, <em>The Tom-Tom Maker</em>. Signed bottom left corner."
$ws.Rows.Item(142).RowHeight = 12.75

# Update description text for 2012.37.23 (row 159, column B)
$ws.Range("B159").Value = "Sculpture by Joe Beeler, <em>Prairie Madonna</em>.  `nInscribed back: <em>JOE BEELER CA</em>`nStamped back: <em>BRONZE/SMITH 18/35.</em>"
$ws.Rows.Item(159).RowHeight = 12.75

# Update description text for 2014.34.1 (row 180, column B) - fix "Billy" to "Bill"
$ws.Range("B180").Value = "Painting by Bill Schenck, <em>A River Runs Through It</em>, 2011."
$ws.Rows.Item(180).RowHeight = 12.75

# Update the active selection on the sheet view to B2
$ws.Range("B2").Select()

# Update the workbook window view settings (position/size of the workbook window)
$win = $excel.ActiveWindow
$win.Top = 465
$win.Height = 9645
